# Applies the "Add files via upload" revision:
#  - Week 2 matches revert from completed/winner-known back to TBD / not completed
#  - Week 2 picks (which depended on those matches) are removed
#  - Team_Winner_Usage rows are recomputed to reflect only the remaining (Week 1) picks
#  - Users and Summary scores drop by 1 (the now-removed Week 2 point)

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Matches sheet: Week 2 games (rows 18-30) go back to "TBD" / not completed
# ------------------------------------------------------------------
$matches = $wb.Worksheets.Item("Matches")
for ($row = 18; $row -le 30; $row++) {
    $matches.Range("F" + $row).Value = "TBD"
    $matches.Range("G" + $row).Value = $false
}

# ------------------------------------------------------------------
# 2) Picks sheet: drop the Week 2 picks (old rows 6-9), leaving only
#    the Week 1 picks in rows 2-5
# ------------------------------------------------------------------
$picks = $wb.Worksheets.Item("Picks")
$picks.Range("A6:I9").EntireRow.Delete()

# ------------------------------------------------------------------
# 3) Team_Winner_Usage sheet: only the Week 1 winners remain (one use
#    each), so the old rows 3-5 are overwritten with what used to be
#    in rows 4/6/8 and the trailing rows 6-9 are removed
# ------------------------------------------------------------------
$winnerUsage = $wb.Worksheets.Item("Team_Winner_Usage")
$winnerUsage.Range("A3").Value = "Haunschi"
$winnerUsage.Range("B3").Value = "Washington Commanders"
$winnerUsage.Range("A4").Value = "Manuel"
$winnerUsage.Range("B4").Value = "Atlanta Falcons"
$winnerUsage.Range("A5").Value = "Raff"
$winnerUsage.Range("B5").Value = "Cincinnati Bengals"
$winnerUsage.Range("A6:F9").EntireRow.Delete()

# ------------------------------------------------------------------
# 4) Users sheet: scores drop by 1 (the Week 2 point is gone)
# ------------------------------------------------------------------
$users = $wb.Worksheets.Item("Users")
$users.Range("C2").Value = 0
$users.Range("C3").Value = 1
$users.Range("C4").Value = 1
$users.Range("C5").Value = 1

# ------------------------------------------------------------------
# 5) Summary sheet: mirrors the Users scores
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("C2").Value = 0
$summary.Range("D2").Value = "0 points"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = "1 points"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = "1 points"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = "1 points"
